$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update B/C (start/end dates) and S/T (report dates) for rows 8 through 36 ---
for ($r = 8; $r -le 36; $r++) {
    $ws.Cells.Item($r, 2).Value = 44470   # column B
    $ws.Cells.Item($r, 3).Value = 44561   # column C
    $ws.Cells.Item($r, 19).Value = 44571  # column S
    $ws.Cells.Item($r, 20).Value = 44571  # column T
}

# --- Update column O (advance figures) for the rows whose totals changed ---
$ws.Cells.Item(9, 15).Value = 1579
$ws.Cells.Item(11, 15).Value = 2163
$ws.Cells.Item(12, 15).Value = 2760
$ws.Cells.Item(17, 15).Value = 2663
$ws.Cells.Item(18, 15).Value = 168
$ws.Cells.Item(19, 15).Value = 360
$ws.Cells.Item(20, 15).Value = 2495
$ws.Cells.Item(21, 15).Value = 44
$ws.Cells.Item(23, 15).Value = 14
$ws.Cells.Item(25, 15).Value = 6
$ws.Cells.Item(26, 15).Value = 3
$ws.Cells.Item(27, 15).Value = 3
$ws.Cells.Item(28, 15).Value = 40
$ws.Cells.Item(29, 15).Value = 3
$ws.Cells.Item(30, 15).Value = 143
$ws.Cells.Item(31, 15).Value = 3
$ws.Cells.Item(32, 15).Value = 171
$ws.Cells.Item(34, 15).Value = 5

# --- D35: was a shared-string label (idx 237), becomes the literal number 6 ---
$ws.Range("D35").Value = 6

# --- Column U (note text): "2do trimestre" -> "3er trimestre" (quarter updated) ---
$textQ3 = "El avance de metas correspondiente a la Matriz 1235-21-04 Educación Superior del ejercicio 2021 el 3er trimestre 2021. Respecto a las metas ajustadas se reportan al cierre del ejercicio 2021"
$textQ3DoubleSpace = "El avance de metas correspondiente a la Matriz 1235-21-04 Educación Superior del ejercicio 2021 el 3er  trimestre 2021. Respecto a las metas ajustadas se reportan al cierre del ejercicio 2021"
$textFormaParte = "El avance de metas correspondiente forma parte de la Matriz  1226-21-01 Fortalecimiento  a la educación en entrega de útiles y uniformes del ejercicio 2021 es del 3er trimestre 2021. Respecto a las metas ajustadas se reportan al cierre del ejercicio 2021"

foreach ($r in 8..34) {
    if ($r -eq 22) {
        $ws.Cells.Item($r, 21).Value = $textQ3DoubleSpace
    } else {
        $ws.Cells.Item($r, 21).Value = $textQ3
    }
}
$ws.Cells.Item(35, 21).Value = $textFormaParte
$ws.Cells.Item(36, 21).Value = $textFormaParte

# --- Selection / scroll position update ---
$ws.Activate()
$ws.Range("U36").Select()

Write-Output "done"
